# "semaine 25/04 + tests perfs"
# Reorganise the last two weeks of the planning:
#   - split "Tests performances + tests utilisateurs" into two separate tasks
#   - add "About box + assembly info" and the UX-improvements task
#   - clear the task for the new week (25/04, 42492) so it can be filled in later

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert one row right after the "Tests performances + tests
#     utilisateurs" line (row 85), then two more rows just above the
#     separator that used to follow "Aspect visuel amélioré" (now row 88).
[void]$ws.Rows(86).Insert()
[void]$ws.Rows("89:90").Insert()

# --- Row 85 used to be "Tests performances + tests utilisateurs"; it is
#     split in two: "Tests performances" (row 86, new) and "Tests
#     utilisateurs" (row 85, keeps the date but switches to the regular
#     in-week formatting since it is no longer the last line of the week).
[void]$ws.Range("A81").Copy()
[void]$ws.Range("A85").PasteSpecial(-4122)

[void]$ws.Range("C82").Copy()
[void]$ws.Range("C86").PasteSpecial(-4122)
$ws.Range("C86").Value = "Tests performances"

[void]$ws.Range("C82").Copy()
[void]$ws.Range("C85").PasteSpecial(-4122)
$ws.Range("C85").Value = "Tests utilisateurs"

# --- Rows 87-88 already hold "Debug sélection arbre" / "Aspect visuel
#     amélioré" (shifted down from 86/87) with the correct formatting.

# --- Rows 89-90 (new): two extra tasks for that week.
[void]$ws.Range("C82").Copy()
[void]$ws.Range("C89").PasteSpecial(-4122)
$ws.Range("C89").Value = "About box + assembly info"

[void]$ws.Range("C82").Copy()
[void]$ws.Range("C90").PasteSpecial(-4122)
$ws.Range("C90").Value = "Améliorations UX (autoscroll lors de navigation depuis recherche, tooltips, …)"

# --- Row 92 (was 89): the task for week 25/04 (42492) is cleared out,
#     ready to be filled in later.
$ws.Range("C92").Value = ""

# --- Restore the on-screen scroll/selection state to what it was left at.
$excel.ActiveWindow.ScrollRow = 73
[void]$ws.Range("C92").Select()

$excel.CutCopyMode = $false
